$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.870.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").Value = "'2.666.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.51%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'598.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.40%  "

$ws.Range("D6").Value = "'158.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.08%  "

$ws.Range("D7").Value = "'0.653"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.92%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "'0.128"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.70%  "

$ws.Range("E10").Value = "  +0.61%  "

$ws.Range("E11").Value = "  -0.62%  "

$ws.Range("E12").Value = "  +1.58%  "

$ws.Range("D13").Value = "'29.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.29%  "

$ws.Range("D14").Value = "'0.0000195"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.11%  "

$ws.Range("D15").Value = "'3.144.84"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.52%  "

$ws.Range("D16").Value = "'65.712.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.12%  "

$ws.Range("D17").Value = "'2.689.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.00%  "

$ws.Range("D18").Value = "'12.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.52%  "

$ws.Range("D19").Value = "'4.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.03%  "

$ws.Range("D20").Value = "'7.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.16%  "

$ws.Range("D21").Value = "'351.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.30%  "

$ws.Range("E22").Value = "  -0.10%  "

$ws.Range("D23").Value = "'69.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.15%  "

$ws.Range("D24").Value = "'1.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +11.18%  "

$ws.Range("D25").Value = "'0.0000113"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.20%  "

$ws.Range("D26").Value = "'9.67"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.32%  "

$ws.Range("E27").Value = "  +1.44%  "

$ws.Range("D28").Value = "'572.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.86%  "

$ws.Range("D29").Value = "'8.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.12%  "

$ws.Range("E31").Value = "  -0.21%  "

$ws.Range("E32").Value = "  +0.80%  "

$ws.Range("E33").Value = "  +3.71%  "

$ws.Range("D34").Value = "'6.73"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.58%  "

$ws.Range("D35").Value = "'5.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.53%  "

$ws.Range("E36").Value = "  -0.20%  "

$ws.Range("E37").Value = "  +0.21%  "

$ws.Range("E38").Value = "  -0.02%  "

$ws.Range("E39").Value = "  +0.57%  "

$ws.Range("D40").Value = "'154.52"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.44%  "

$ws.Range("D41").Value = "'162.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.87%  "

$ws.Range("D42").Value = "'4.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.88%  "

$ws.Range("D43").Value = "'0.0621"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.50%  "

$ws.Range("E44").Value = "  -0.21%  "

$ws.Range("D45").Value = "'23.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.95%  "

$ws.Range("E46").Value = "  +0.25%  "

$ws.Range("D47").Value = "'0.0259"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.30%  "

$ws.Range("E48").Value = "  +2.12%  "

$ws.Range("D49").Value = "'19.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.34%  "

$ws.Range("E50").Value = "  -5.46%  "

$ws.Range("D51").Value = "'0.818"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.20%  "
